# Balance enquiry, mini statement, delete account and customised statement
# Adds a new "AccountIDForDelete" / "78178" row to the TestData sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 14: label in A14, quote-prefixed numeric-looking id in B14
# (leading apostrophe forces it to be stored as text, same convention
# already used for CustomerID / InitialDeposit above it).
$ws.Range("A14").Value = "AccountIDForDelete"
$ws.Range("B14").Value = "'78178"

# Leave the same cell selected/active as in the authored workbook.
$null = $ws.Range("G13").Select()
